# This script applies the "456a3b4" content update to the 江西-漫展信息 workbook.
# Two sheets ("展览" and "全部类型") each had their first two data rows
# (the two cancelled events: 上饶·次元重现夏日嘉年华（取消） and
#  乐平·CY境界次元第三届动漫游戏庆典（取消）) removed, which shifts every
# subsequent row up by two. After the shift, the running index in column A
# is renumbered, and the "想去人数" (column F) counter is refreshed to the
# latest scraped values.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (dimension becomes A1:I15) ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Rows.Item(2).Delete() | Out-Null
$ws1.Rows.Item(2).Delete() | Out-Null

$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(9, 1).Value = 8
$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(11, 1).Value = 10
$ws1.Cells.Item(12, 1).Value = 11
$ws1.Cells.Item(13, 1).Value = 12
$ws1.Cells.Item(14, 1).Value = 13
$ws1.Cells.Item(15, 1).Value = 14

$ws1.Cells.Item(2, 6).Value = 118
$ws1.Cells.Item(3, 6).Value = 5590
$ws1.Cells.Item(4, 6).Value = 78
$ws1.Cells.Item(5, 6).Value = 5
$ws1.Cells.Item(6, 6).Value = 923
$ws1.Cells.Item(7, 6).Value = 148
$ws1.Cells.Item(8, 6).Value = 2502
$ws1.Cells.Item(9, 6).Value = 83
$ws1.Cells.Item(10, 6).Value = 124
$ws1.Cells.Item(11, 6).Value = 4
$ws1.Cells.Item(12, 6).Value = 75
$ws1.Cells.Item(13, 6).Value = 13
$ws1.Cells.Item(14, 6).Value = 2338
$ws1.Cells.Item(15, 6).Value = 305

# ---- Sheet "全部类型" (dimension becomes A1:I18) ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Rows.Item(2).Delete() | Out-Null
$ws4.Rows.Item(2).Delete() | Out-Null

$ws4.Cells.Item(2, 1).Value = 1
$ws4.Cells.Item(3, 1).Value = 2
$ws4.Cells.Item(4, 1).Value = 3
$ws4.Cells.Item(5, 1).Value = 4
$ws4.Cells.Item(6, 1).Value = 5
$ws4.Cells.Item(7, 1).Value = 6
$ws4.Cells.Item(8, 1).Value = 7
$ws4.Cells.Item(9, 1).Value = 8
$ws4.Cells.Item(10, 1).Value = 9
$ws4.Cells.Item(11, 1).Value = 10
$ws4.Cells.Item(12, 1).Value = 11
$ws4.Cells.Item(13, 1).Value = 12
$ws4.Cells.Item(14, 1).Value = 13
$ws4.Cells.Item(15, 1).Value = 14
$ws4.Cells.Item(16, 1).Value = 15
$ws4.Cells.Item(17, 1).Value = 16
$ws4.Cells.Item(18, 1).Value = 17

$ws4.Cells.Item(2, 6).Value = 118
$ws4.Cells.Item(3, 6).Value = 5590
$ws4.Cells.Item(4, 6).Value = 101
$ws4.Cells.Item(5, 6).Value = 78
$ws4.Cells.Item(6, 6).Value = 5
$ws4.Cells.Item(7, 6).Value = 2
$ws4.Cells.Item(8, 6).Value = 923
$ws4.Cells.Item(9, 6).Value = 148
$ws4.Cells.Item(10, 6).Value = 2502
$ws4.Cells.Item(11, 6).Value = 83
$ws4.Cells.Item(12, 6).Value = 124
$ws4.Cells.Item(13, 6).Value = 4
$ws4.Cells.Item(14, 6).Value = 1
$ws4.Cells.Item(15, 6).Value = 75
$ws4.Cells.Item(16, 6).Value = 13
$ws4.Cells.Item(17, 6).Value = 2338
$ws4.Cells.Item(18, 6).Value = 305
